# moved sampling parameters from 2EXT_01 to 1SPL01
# Remove the "Sample Collection Method", "Metabolism quenching method" and
# "Sample storage" parameter column groups (9 columns total: Parameter +
# Term Source REF + Term Accession Number, x3) from the annotation table.
# These used to live in columns C:K; deleting them shifts the remaining
# "Bio entity" / "Biosource amount" / "Extraction Kit" columns left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Remember table cosmetic properties before we rebuild it.
$tblName = $tbl.Name
$tblStyle = "TableStyleMedium7"

# Physically remove the 9 obsolete columns (C through K), which shifts
# everything to the right of them (Bio entity onward) left into place.
$ws.Range("C1:K1").EntireColumn.Delete()

# The table definition itself does not auto-shrink when interior
# worksheet columns are deleted, so rebuild it against the new extent
# to get column headers/count back in sync with the sheet data.
$tbl.Unlist()
$newTbl = $ws.ListObjects.Add(1, $ws.Range("A2:N3"), $null, 1)
$newTbl.Name = $tblName
$newTbl.TableStyle = $tblStyle

# Update the active selection to match the post-edit workbook state.
$ws.Range("T12").Select()
